$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: environment changed from the "i-preproduccion" QA box to the
# "ssurgwsoadev4-oci" Oracle-hosted gateway box (PC Gestion Documental -> R).
#   A2 Ambiente     -> ssurgwsoadev4-oci.opc.oracleoutsourcing.com
#   B2 URL          -> https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do (+ hyperlink)
#   C2 Usuario      -> su
#   D2 Contrasenia  -> gw
#   E2 NroPoliza    -> " 04104016708" (leading space kept, stored as text)

# A2: use a leading apostrophe so the engine keeps it as text and keeps the
# cell's existing (quote-prefixed) style instead of resetting it.
$ws.Range("A2").Formula = "'ssurgwsoadev4-oci.opc.oracleoutsourcing.com"

# C2 / D2 are plain, unstyled text cells.
$ws.Range("C2").Formula = "su"
$ws.Range("D2").Formula = "gw"

# E2 keeps a leading space and must stay text (not get coerced to a number) -
# the apostrophe prefix forces text storage while preserving the cell style.
$ws.Range("E2").Formula = "' 04104016708"

# B2: replace the hyperlink target + display text. Remove the old hyperlink
# first (this cleanly drops the old relationship), set the new text, re-add
# the hyperlink, then restore the "Hipervínculo" cell style that Add() resets.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B2").Formula = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do")
$ws.Range("B2").Style = "Hipervínculo"

# Selection moves to A2:D2 with A2 active, matching the saved view state.
[void]$ws.Range("A2:D2").Select()
